$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$values = @{
    "H2" = 2304.0715
    "I2" = 506.33334
    "J2" = 2794.3635
    "K2" = 506.33334
    "L2" = 2794.3635
    "M2" = -393.33334
    "N2" = -3020.3635
    "H21" = 49900
    "I21" = 0
    "J21" = 49900
    "K21" = 0
    "L21" = 49900
    "N21" = -50836
    "H23" = 49900
    "I23" = 0
    "J23" = 49900
    "K23" = 0
    "L23" = 49900
    "N23" = -50368
    "H28" = 437.18182
    "I28" = 201.1
    "K28" = 201.1
    "M28" = 283.9
    "H29" = 0
    "I29" = 0
    "J29" = 0
    "K29" = 0
    "L29" = 0
    "H38" = 1686.1
    "I38" = 1686.1
    "K38" = 5058.299999999999
    "M38" = -4686.299999999999
    "H40" = 7439
    "I40" = 7050
    "J40" = 8995
    "K40" = 7050
    "L40" = 8995
    "M40" = -6875
    "N40" = -9345
    "H43" = 11406.363
    "I43" = 1978.6
    "J43" = 31608.715
    "K43" = 1978.6
    "L43" = 31608.715
    "M43" = -1909.6
    "N43" = -31746.715
    "H53" = 165.2
    "I53" = 63.4
    "K53" = 63.4
    "M53" = 573.6
    "H55" = 677.1667
    "I55" = 381
    "K55" = 381
    "M55" = -167
    "H92" = 1552.5555
    "I92" = 530.06665
    "K92" = 530.06665
    "M92" = 717.93335
    "H115" = 751.4545
    "I115" = 745.1111
    "J115" = 780
    "K115" = 2235.3333
    "L115" = 2340
    "M115" = -668.3332999999998
    "N115" = -5474
    "H121" = 3482.5
    "J121" = 3482.5
    "L121" = 10447.5
    "N121" = -13941.5
    "H129" = 1855.2
    "I129" = 789
    "J129" = 2566
    "K129" = 2367
    "L129" = 7698
    "M129" = 2633
    "N129" = -17698
    "H132" = 436816.6
    "I132" = 2104.4614
    "J132" = 2858784.2
    "K132" = 6313.3842
    "L132" = 8576352.600000001
    "M132" = -3783.3842
    "N132" = -8581412.600000001
    "H135" = 1498
    "I135" = 1555.3846
    "J135" = 1125
    "K135" = 13998.4614
    "L135" = 10125
    "M135" = -11463.4614
    "N135" = -15195
    "H138" = 3086.074
    "I138" = 2978.9167
    "J138" = 3116.6904
    "K138" = 8936.750100000001
    "L138" = 9350.0712
    "M138" = -3796.750100000001
    "N138" = -19630.0712
}
foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
$toClear = @("M21", "M23", "M29", "N29")
foreach ($addr in $toClear) {
    $ws.Range($addr).ClearContents()
}

$ws = $wb.Worksheets.Item("ARM")
$values = @{
    "H2" = 1552.2
    "J2" = 1806.5714
    "L2" = 1806.5714
    "N2" = -2032.5714
    "H32" = 21879.176
    "I32" = 8344.058
    "K32" = 8344.058
    "M32" = -8057.058000000001
    "H45" = 1395.6
    "I45" = 1495
    "K45" = 1495
    "M45" = -1118
    "H61" = 3483.1428
    "I61" = 2419.75
    "K61" = 2419.75
    "M61" = -2207.75
    "H74" = 956.4
    "I74" = 956.4
    "K74" = 956.4
    "M74" = -82.39999999999998
    "H77" = 956.4
    "I77" = 956.4
    "K77" = 4782
    "M77" = -414
    "H97" = 1108.129
    "I97" = 1068.0416
    "K97" = 1068.0416
    "M97" = -572.0416
    "H102" = 2227.6667
    "I102" = 2049.913
    "J102" = 3249.75
    "K102" = 2049.913
    "L102" = 3249.75
    "M102" = -427.913
    "N102" = -6493.75
    "H110" = 2034.4166
    "I110" = 1626.4
    "K110" = 1626.4
    "M110" = 418.5999999999999
    "H116" = 1552.2
    "J116" = 1806.5714
    "L116" = 1806.5714
    "N116" = -6394.5714
    "H122" = 2919.325
    "I122" = 2843
    "J122" = 3224.625
    "K122" = 8529
    "L122" = 9673.875
    "M122" = -6079
    "N122" = -14573.875
    "H132" = 2083.5264
    "I132" = 1009.93335
    "J132" = 6109.5
    "K132" = 3029.80005
    "L132" = 18328.5
    "M132" = -499.8000499999998
    "N132" = -23388.5
    "H136" = 3483.1428
    "I136" = 2419.75
    "K136" = 7259.25
    "M136" = -4709.25
    "H138" = 62000
    "J138" = 62000
    "L138" = 62000
    "N138" = -72280
}
foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

$ws = $wb.Worksheets.Item("BSM")
$values = @{
    "H3" = 1552.2
    "J3" = 1806.5714
    "L3" = 1806.5714
    "N3" = -2034.5714
    "H22" = 584.5909
    "J22" = 735.3333
    "L22" = 735.3333
    "N22" = -1081.3333
    "H29" = 7166.6665
    "I29" = 7166.6665
    "K29" = 7166.6665
    "M29" = -6877.6665
    "H81" = 66374.25
    "J81" = 71833
    "L81" = 71833
    "N81" = -73955
    "H84" = 66374.25
    "J84" = 71833
    "L84" = 215499
    "N84" = -226107
    "H94" = 2681.577
    "I94" = 1831.35
    "K94" = 1831.35
    "M94" = -1380.35
    "H99" = 1164.1111
    "I99" = 1110.9231
    "J99" = 1302.4
    "K99" = 1110.9231
    "L99" = 1302.4
    "M99" = 387.0769
    "N99" = -4298.4
    "H105" = 2882.2222
    "I105" = 2248.1428
    "J105" = 3565.077
    "K105" = 2248.1428
    "L105" = 3565.077
    "M105" = -501.1428000000001
    "N105" = -7059.077
    "H134" = 2756.4
    "I134" = 2600.4167
    "J134" = 6500
    "K134" = 7801.250100000001
    "L134" = 19500
    "M134" = -5266.250100000001
    "N134" = -24570
    "H138" = 84833.336
    "J138" = 84833.336
    "L138" = 84833.336
    "N138" = -95113.336
    "H139" = 91000
    "J139" = 82000
    "L139" = 82000
    "N139" = -92280
}
foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

$ws = $wb.Worksheets.Item("CRP")
$values = @{
    "H7" = 67119.07
    "I7" = 111234.89
    "K7" = 111234.89
    "M7" = -111121.89
    "H16" = 1372.3636
    "J16" = 1185.25
    "L16" = 1185.25
    "N16" = -1759.25
    "H58" = 4056.4707
    "I58" = 4243.1333
    "K58" = 4243.1333
    "M58" = -4040.1333
    "H86" = 8506
    "I86" = 5924
    "J86" = 9366.667
    "K86" = 5924
    "L86" = 9366.667
    "M86" = -4801
    "N86" = -11612.667
    "H89" = 8506
    "I89" = 5924
    "J89" = 9366.667
    "K89" = 29620
    "L89" = 46833.335
    "M89" = -24004
    "N89" = -58065.335
    "H113" = 1372.3636
    "J113" = 1185.25
    "L113" = 1185.25
    "N113" = -5525.25
    "H132" = 2446.818
    "I132" = 2624.5
    "K132" = 7873.5
    "M132" = -5343.5
    "H134" = 3302.25
    "I134" = 3302.25
    "J134" = 0
    "K134" = 9906.75
    "L134" = 0
    "M134" = -7371.75
    "H136" = 4056.4707
    "I136" = 4243.1333
    "K136" = 12729.3999
    "M136" = -10179.3999
    "H141" = 556499.25
    "J141" = 556499.25
    "L141" = 556499.25
    "N141" = -566859.25
}
foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
$toClear = @("N134")
foreach ($addr in $toClear) {
    $ws.Range($addr).ClearContents()
}

$ws = $wb.Worksheets.Item("CUL")
$values = @{
    "H2" = 410.75
    "I2" = 150
    "J2" = 497.66666
    "K2" = 900
    "L2" = 2985.99996
    "M2" = -787
    "N2" = -3211.99996
    "H3" = 6550.6665
    "I3" = 2986.6667
    "J3" = 8332.667
    "K3" = 8960.000100000001
    "L3" = 24998.001
    "M3" = -8848.000100000001
    "N3" = -25222.001
    "H5" = 696.8
    "J5" = 993.5
    "L5" = 2980.5
    "N5" = -3204.5
    "H38" = 124.21429
    "I38" = 135.44444
    "J38" = 104
    "K38" = 406.33332
    "L38" = 312
    "M38" = -59.33331999999996
    "N38" = -1006
    "H70" = 966.3333
    "I70" = 966.3333
    "K70" = 2898.9999
    "M70" = -2583.9999
    "H73" = 966.3333
    "I73" = 966.3333
    "K73" = 2898.9999
    "M73" = -1806.9999
    "H113" = 794.7619
    "I113" = 398.25
    "J113" = 888.05884
    "K113" = 1194.75
    "L113" = 2664.17652
    "M113" = 975.25
    "N113" = -7004.17652
    "H135" = 696.8
    "J135" = 993.5
    "L135" = 8941.5
    "N135" = -14011.5
}
foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

$ws = $wb.Worksheets.Item("GSM")
$values = @{
    "H3" = 1738.1428
    "I3" = 943
    "J3" = 1799.3077
    "K3" = 943
    "L3" = 1799.3077
    "M3" = -827
    "N3" = -2031.3077
    "H5" = 7600.1333
    "I5" = 8615.692
    "J5" = 999
    "K5" = 8615.692
    "L5" = 999
    "M5" = -8503.692
    "N5" = -1223
    "H26" = 45018.5
    "J26" = 0
    "L26" = 0
    "H50" = 45018.5
    "J50" = 0
    "L50" = 0
    "H80" = 9275.083
    "I80" = 2117.5
    "J80" = 16432.666
    "K80" = 2117.5
    "L80" = 16432.666
    "M80" = -1119.5
    "N80" = -18428.666
    "H83" = 9275.083
    "I83" = 2117.5
    "J83" = 16432.666
    "K83" = 10587.5
    "L83" = 82163.33
    "M83" = -5595.5
    "N83" = -92147.33
    "H95" = 35172
    "J95" = 35172
    "L95" = 35172
    "N95" = -40664
    "H102" = 2924.5
    "I102" = 2679.0667
    "K102" = 2679.0667
    "M102" = -1057.0667
    "H122" = 1920.3846
    "I122" = 1503.5
    "J122" = 1996.1818
    "K122" = 4510.5
    "L122" = 5988.5454
    "M122" = -2060.5
    "N122" = -10888.5454
    "H126" = 2940.7058
    "I126" = 2777
    "J126" = 3124.875
    "K126" = 8331
    "L126" = 9374.625
    "M126" = -5861
    "N126" = -14314.625
    "H136" = 15218.4
    "J136" = 15218.4
    "L136" = 45655.2
    "N136" = -50755.2
}
foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
$toClear = @("N26", "N50")
foreach ($addr in $toClear) {
    $ws.Range($addr).ClearContents()
}

$ws = $wb.Worksheets.Item("LTW")
$values = @{
    "H7" = 5972.364
    "I7" = 5309.619
    "J7" = 19890
    "K7" = 5309.619
    "L7" = 19890
    "M7" = -5197.619
    "N7" = -20114
    "H35" = 1500
    "I35" = 1500
    "J35" = 0
    "K35" = 1500
    "L35" = 0
    "M35" = -1164
    "H82" = 2109.7144
    "I82" = 1962.4166
    "J82" = 2993.5
    "K82" = 1962.4166
    "L82" = 2993.5
    "M82" = -1601.4166
    "N82" = -3715.5
    "H85" = 2109.7144
    "I85" = 1962.4166
    "J85" = 2993.5
    "K85" = 1962.4166
    "L85" = 2993.5
    "M85" = -714.4166
    "N85" = -5489.5
    "H122" = 2931.5715
    "I122" = 2452
    "K122" = 7356
    "M122" = -4906
    "H126" = 5972.364
    "I126" = 5309.619
    "J126" = 19890
    "K126" = 15928.857
    "L126" = 59670
    "M126" = -13458.857
    "N126" = -64610
    "H132" = 3091.6843
    "I132" = 3249.5715
    "K132" = 9748.7145
    "M132" = -7218.7145
    "H134" = 59214
    "J134" = 59214
    "L134" = 59214
    "N134" = -69354
}
foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
$toClear = @("N35")
foreach ($addr in $toClear) {
    $ws.Range($addr).ClearContents()
}

$ws = $wb.Worksheets.Item("WVR")
$values = @{
    "H3" = 894333.3
    "I3" = 1067000
    "J3" = 31000
    "K3" = 1067000
    "L3" = 31000
    "M3" = -1066886
    "N3" = -31228
    "H12" = 6
    "I12" = 6
    "K12" = 6
    "M12" = 136
    "H33" = 20021
    "J33" = 20021
    "L33" = 20021
    "N33" = -20521
    "H36" = 20021
    "J36" = 20021
    "L36" = 20021
    "N36" = -20521
    "H39" = 26283
    "I39" = 22500
    "J39" = 27796.2
    "K39" = 22500
    "L39" = 27796.2
    "M39" = -22087
    "N39" = -28622.2
    "H42" = 49999
    "I42" = 49999
    "K42" = 49999
    "M42" = -49621
    "H81" = 2198.3333
    "I81" = 780
    "J81" = 3616.6667
    "K81" = 1560
    "L81" = 7233.3334
    "M81" = -499
    "N81" = -9355.3334
    "H84" = 2198.3333
    "I84" = 780
    "J84" = 3616.6667
    "K84" = 7800
    "L84" = 36166.667
    "M84" = -2496
    "N84" = -46774.667
    "H107" = 1005.4
    "I107" = 671.7778
    "K107" = 2015.3334
    "M107" = -95.33339999999998
    "H122" = 10177.226
    "I122" = 9880.929
    "K122" = 29642.787
    "M122" = -27192.787
    "H126" = 3324.8076
    "I126" = 2521.2856
    "J126" = 6699.6
    "K126" = 7563.8568
    "L126" = 20098.8
    "M126" = -5093.8568
    "N126" = -25038.8
    "H132" = 1942.5333
    "I132" = 2119.36
    "J132" = 1058.4
    "K132" = 6358.08
    "L132" = 3175.2
    "M132" = -3828.08
    "N132" = -8235.2
    "H137" = 90000
    "J137" = 90000
    "L137" = 90000
    "N137" = -100200
    "H141" = 0
    "J141" = 0
    "L141" = 0
}
foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
$toClear = @("N141")
foreach ($addr in $toClear) {
    $ws.Range($addr).ClearContents()
}
